$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.934
$ws.Range("D6").Value = -7.766
$ws.Range("D7").Value = -7.258999999999999
$ws.Range("C8").Value = -12.672
$ws.Range("D8").Value = -7.858
$ws.Range("E11").Value = 12.681
$ws.Range("B12").Value = 5.513
$ws.Range("C12").Value = -13.073
$ws.Range("C14").Value = -11.843
$ws.Range("E14").Value = 12.915
$ws.Range("D19").Value = -7.74
$ws.Range("E19").Value = 12.67
$ws.Range("D21").Value = -7.478999999999999
$ws.Range("E21").Value = 13.378
$ws.Range("C22").Value = -12.601
$ws.Range("D24").Value = -7.934
